$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 590.4
$ws.Range("I55").Value = 376.25
$ws.Range("J55").Value = 733.1667
$ws.Range("K55").Value = 376.25
$ws.Range("L55").Value = 733.1667
$ws.Range("M55").Value = -162.25
$ws.Range("N55").Value = -1161.1667

$ws.Range("H88").Value = 1125197.8
$ws.Range("I88").Value = 753
$ws.Range("J88").Value = 1375074.4
$ws.Range("K88").Value = 753
$ws.Range("L88").Value = 1375074.4
$ws.Range("M88").Value = -347
$ws.Range("N88").Value = -1375886.4

$ws.Range("H91").Value = 1125197.8
$ws.Range("I91").Value = 753
$ws.Range("J91").Value = 1375074.4
$ws.Range("K91").Value = 753
$ws.Range("L91").Value = 1375074.4
$ws.Range("M91").Value = 651
$ws.Range("N91").Value = -1377882.4

$ws.Range("H138").Value = 513230.16
$ws.Range("I138").Value = 1798.5
$ws.Range("J138").Value = 581421.0600000001
$ws.Range("K138").Value = 5395.5
$ws.Range("L138").Value = 1744263.18
$ws.Range("M138").Value = -255.5
$ws.Range("N138").Value = -1754543.18


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1085.7142
$ws.Range("I2").Value = 900
$ws.Range("J2").Value = 1225
$ws.Range("K2").Value = 900
$ws.Range("L2").Value = 1225
$ws.Range("M2").Value = -787
$ws.Range("N2").Value = -1451

$ws.Range("H45").Value = 1507
$ws.Range("I45").Value = 1396.6666
$ws.Range("K45").Value = 1396.6666
$ws.Range("M45").Value = -1019.6666

$ws.Range("H116").Value = 1085.7142
$ws.Range("I116").Value = 900
$ws.Range("J116").Value = 1225
$ws.Range("K116").Value = 900
$ws.Range("L116").Value = 1225
$ws.Range("M116").Value = 1394
$ws.Range("N116").Value = -5813

$ws.Range("H132").Value = 3788.5715
$ws.Range("I132").Value = 3616
$ws.Range("K132").Value = 10848
$ws.Range("M132").Value = -8318


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1085.7142
$ws.Range("I3").Value = 900
$ws.Range("J3").Value = 1225
$ws.Range("K3").Value = 900
$ws.Range("L3").Value = 1225
$ws.Range("M3").Value = -786
$ws.Range("N3").Value = -1453

$ws.Range("H13").Value = 24999
$ws.Range("J13").Value = 24999
$ws.Range("L13").Value = 24999
$ws.Range("N13").Value = -25335

$ws.Range("H94").Value = 20833902
$ws.Range("I94").Value = 22727802
$ws.Range("K94").Value = 22727802
$ws.Range("M94").Value = -22727351

$ws.Range("H107").Value = 1231.5883
$ws.Range("I107").Value = 856.7273
$ws.Range("K107").Value = 856.7273
$ws.Range("M107").Value = 1063.2727

$ws.Range("H109").Value = 10000
$ws.Range("J109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("N109").Value = -12774


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1423.8667
$ws.Range("I31").Value = 1489.5264
$ws.Range("J31").Value = 1375.8846
$ws.Range("K31").Value = 1489.5264
$ws.Range("L31").Value = 1375.8846
$ws.Range("M31").Value = -1194.5264
$ws.Range("N31").Value = -1965.8846

$ws.Range("H34").Value = 1423.8667
$ws.Range("I34").Value = 1489.5264
$ws.Range("J34").Value = 1375.8846
$ws.Range("K34").Value = 1489.5264
$ws.Range("L34").Value = 1375.8846
$ws.Range("M34").Value = -1287.5264
$ws.Range("N34").Value = -1779.8846

$ws.Range("H81").Value = 15000
$ws.Range("J81").Value = 15000
$ws.Range("L81").Value = 15000
$ws.Range("N81").Value = -16996

$ws.Range("H84").Value = 15000
$ws.Range("J84").Value = 15000
$ws.Range("L84").Value = 45000
$ws.Range("N84").Value = -54984

$ws.Range("H99").Value = 1517
$ws.Range("J99").Value = 1415.2
$ws.Range("L99").Value = 1415.2
$ws.Range("N99").Value = -4411.2

$ws.Range("H107").Value = 940.93335
$ws.Range("I107").Value = 483.0909
$ws.Range("J107").Value = 2200
$ws.Range("K107").Value = 483.0909
$ws.Range("L107").Value = 2200
$ws.Range("M107").Value = 1436.9091
$ws.Range("N107").Value = -6040

$ws.Range("H114").Value = 23998
$ws.Range("J114").Value = 23998
$ws.Range("L114").Value = 23998
$ws.Range("N114").Value = -32676

$ws.Range("H126").Value = 1517
$ws.Range("J126").Value = 1415.2
$ws.Range("L126").Value = 4245.6
$ws.Range("N126").Value = -9185.6

$ws.Range("H133").Value = 28054.285
$ws.Range("J133").Value = 28054.285
$ws.Range("L133").Value = 28054.285
$ws.Range("N133").Value = -33114.285


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10068.228
$ws.Range("I3").Value = 3570.077
$ws.Range("K3").Value = 10710.231
$ws.Range("M3").Value = -10598.231

$ws.Range("H4").Value = 1758157
$ws.Range("I4").Value = 3334366.2
$ws.Range("J4").Value = 576000
$ws.Range("K4").Value = 10003098.6
$ws.Range("L4").Value = 1728000
$ws.Range("M4").Value = -10002986.6
$ws.Range("N4").Value = -1728224

$ws.Range("H70").Value = 6419.9
$ws.Range("I70").Value = 3999
$ws.Range("J70").Value = 6688.8887
$ws.Range("K70").Value = 11997
$ws.Range("L70").Value = 20066.6661
$ws.Range("M70").Value = -11682
$ws.Range("N70").Value = -20696.6661

$ws.Range("H73").Value = 6419.9
$ws.Range("I73").Value = 3999
$ws.Range("J73").Value = 6688.8887
$ws.Range("K73").Value = 11997
$ws.Range("L73").Value = 20066.6661
$ws.Range("M73").Value = -10905
$ws.Range("N73").Value = -22250.6661

$ws.Range("H74").Value = 5750
$ws.Range("J74").Value = 5750
$ws.Range("L74").Value = 17250
$ws.Range("N74").Value = -19372

$ws.Range("H77").Value = 5750
$ws.Range("J77").Value = 5750
$ws.Range("L77").Value = 51750
$ws.Range("N77").Value = -62358

$ws.Range("H87").Value = 2838
$ws.Range("I87").Value = 1014
$ws.Range("J87").Value = 3750
$ws.Range("K87").Value = 3042
$ws.Range("L87").Value = 11250
$ws.Range("M87").Value = -1794
$ws.Range("N87").Value = -13746

$ws.Range("H88").Value = 6084.615
$ws.Range("J88").Value = 7045.4546
$ws.Range("L88").Value = 21136.3638
$ws.Range("N88").Value = -21992.3638

$ws.Range("H90").Value = 2838
$ws.Range("I90").Value = 1014
$ws.Range("J90").Value = 3750
$ws.Range("K90").Value = 9126
$ws.Range("L90").Value = 33750
$ws.Range("M90").Value = -2886
$ws.Range("N90").Value = -46230

$ws.Range("H91").Value = 6084.615
$ws.Range("J91").Value = 7045.4546
$ws.Range("L91").Value = 21136.3638
$ws.Range("N91").Value = -24100.3638

$ws.Range("H132").Value = 803.7619
$ws.Range("I132").Value = 810.3077
$ws.Range("J132").Value = 793.125
$ws.Range("K132").Value = 7292.7693
$ws.Range("L132").Value = 7138.125
$ws.Range("M132").Value = -4762.7693
$ws.Range("N132").Value = -12198.125


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1083.3158
$ws.Range("I102").Value = 989.1539
$ws.Range("J102").Value = 1287.3334
$ws.Range("K102").Value = 989.1539
$ws.Range("L102").Value = 1287.3334
$ws.Range("M102").Value = 632.8461
$ws.Range("N102").Value = -4531.3334

$ws.Range("H126").Value = 1930.2222
$ws.Range("I126").Value = 1624.5714
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4873.7142
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -2403.7142
$ws.Range("N126").Value = -13940


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1726.25
$ws.Range("J7").Value = 1635
$ws.Range("L7").Value = 1635
$ws.Range("N7").Value = -1859

$ws.Range("H22").Value = 880.3333
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 956.4
$ws.Range("K22").Value = 500
$ws.Range("L22").Value = 956.4
$ws.Range("M22").Value = -205
$ws.Range("N22").Value = -1546.4

$ws.Range("H27").Value = 880.3333
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 956.4
$ws.Range("K27").Value = 500
$ws.Range("L27").Value = 956.4
$ws.Range("M27").Value = -393
$ws.Range("N27").Value = -1170.4

$ws.Range("H40").Value = 3326.25
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 3652.5
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 3652.5
$ws.Range("M40").Value = -2864
$ws.Range("N40").Value = -3924.5

$ws.Range("H55").Value = 253.29411
$ws.Range("I55").Value = 160.33333
$ws.Range("J55").Value = 476.4
$ws.Range("K55").Value = 160.33333
$ws.Range("L55").Value = 476.4
$ws.Range("M55").Value = 12.66667000000001
$ws.Range("N55").Value = -822.4

$ws.Range("H68").Value = 1800.4
$ws.Range("I68").Value = 1800.4
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1800.4
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -1051.4
$ws.Range("N68").Value = $null

$ws.Range("H71").Value = 1800.4
$ws.Range("I71").Value = 1800.4
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 9002
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -5258
$ws.Range("N71").Value = $null

$ws.Range("H122").Value = 27779800
$ws.Range("I122").Value = 41668500
$ws.Range("K122").Value = 125005500
$ws.Range("M122").Value = -125003050

$ws.Range("H123").Value = 40950
$ws.Range("J123").Value = 40950
$ws.Range("L123").Value = 40950
$ws.Range("N123").Value = -50750

$ws.Range("H126").Value = 1726.25
$ws.Range("J126").Value = 1635
$ws.Range("L126").Value = 4905
$ws.Range("N126").Value = -9845

$ws.Range("H130").Value = 44809.668
$ws.Range("J130").Value = 44809.668
$ws.Range("L130").Value = 44809.668
$ws.Range("N130").Value = -54849.668


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 52632390
$ws.Range("I126").Value = 62500600
$ws.Range("J126").Value = 1966.3334
$ws.Range("K126").Value = 187501800
$ws.Range("L126").Value = 5899.0002
$ws.Range("M126").Value = -187499330
$ws.Range("N126").Value = -10839.0002

$ws.Range("H132").Value = 1764.8379
$ws.Range("I132").Value = 1525.0333
$ws.Range("J132").Value = 2792.5715
$ws.Range("K132").Value = 4575.0999
$ws.Range("L132").Value = 8377.7145
$ws.Range("M132").Value = -2045.0999
$ws.Range("N132").Value = -13437.7145

$ws.Range("H136").Value = 1215.3462
$ws.Range("I136").Value = 1185.9524
$ws.Range("J136").Value = 1338.8
$ws.Range("K136").Value = 3557.857199999999
$ws.Range("L136").Value = 4016.4
$ws.Range("M136").Value = -1007.857199999999
$ws.Range("N136").Value = -9116.4

